# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 2
    7  = 3
    8  = 9
    9  = 3
    10 = 3
    11 = 3
    12 = 1
    13 = 4
    14 = 6
    15 = 0
    16 = 2
    17 = 3
    18 = 2
    19 = 2
    20 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
